$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark attendance ("1") for the 23/2 column (column G) on rows 4, 6 and 7,
# matching the same value already used for the other attendance marks.
$ws.Range("G4").Value = "1"
$ws.Range("G6").Value = "1"
$ws.Range("G7").Value = "1"

# Reflect the new active selection left by the edit.
$ws.Range("G7").Select()
